$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 2766.6667
$ws.Cells.Item(125, 9).Value = 1300
$ws.Cells.Item(125, 10).Value = 3500
$ws.Cells.Item(125, 11).Value = 11700
$ws.Cells.Item(125, 12).Value = 31500
$ws.Cells.Item(125, 13).Value = -9240
$ws.Cells.Item(125, 14).Value = -36420
$ws.Cells.Item(127, 8).Value = 1611.1666
$ws.Cells.Item(127, 9).Value = 566.6667
$ws.Cells.Item(127, 10).Value = 1959.3334
$ws.Cells.Item(127, 11).Value = 1700.0001
$ws.Cells.Item(127, 12).Value = 5878.0002
$ws.Cells.Item(127, 13).Value = 3259.9999
$ws.Cells.Item(127, 14).Value = -15798.0002
$ws.Cells.Item(129, 8).Value = 3522579.2
$ws.Cells.Item(129, 9).Value = 31251044
$ws.Cells.Item(129, 10).Value = 1504.2222
$ws.Cells.Item(129, 11).Value = 93753132
$ws.Cells.Item(129, 12).Value = 4512.6666
$ws.Cells.Item(129, 13).Value = -93748132
$ws.Cells.Item(129, 14).Value = -14512.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1636.0526
$ws.Cells.Item(45, 9).Value = 1073.8334
$ws.Cells.Item(45, 11).Value = 1073.8334
$ws.Cells.Item(45, 13).Value = -696.8334
$ws.Cells.Item(52, 8).Value = 39779.617
$ws.Cells.Item(52, 10).Value = 39779.617
$ws.Cells.Item(52, 12).Value = 39779.617
$ws.Cells.Item(52, 14).Value = -40415.617
$ws.Cells.Item(61, 8).Value = 3783.7144
$ws.Cells.Item(61, 9).Value = 1746.5
$ws.Cells.Item(61, 10).Value = 6500
$ws.Cells.Item(61, 11).Value = 1746.5
$ws.Cells.Item(61, 12).Value = 6500
$ws.Cells.Item(61, 13).Value = -1534.5
$ws.Cells.Item(61, 14).Value = -6924
$ws.Cells.Item(74, 8).Value = 1542
$ws.Cells.Item(74, 9).Value = 1486.8334
$ws.Cells.Item(74, 10).Value = 1624.75
$ws.Cells.Item(74, 11).Value = 1486.8334
$ws.Cells.Item(74, 12).Value = 1624.75
$ws.Cells.Item(74, 13).Value = -612.8334
$ws.Cells.Item(74, 14).Value = -3372.75
$ws.Cells.Item(77, 8).Value = 1542
$ws.Cells.Item(77, 9).Value = 1486.8334
$ws.Cells.Item(77, 10).Value = 1624.75
$ws.Cells.Item(77, 11).Value = 7434.166999999999
$ws.Cells.Item(77, 12).Value = 8123.75
$ws.Cells.Item(77, 13).Value = -3066.166999999999
$ws.Cells.Item(77, 14).Value = -16859.75
$ws.Cells.Item(122, 8).Value = 3091.0557
$ws.Cells.Item(122, 9).Value = 1969.9166
$ws.Cells.Item(122, 10).Value = 5333.3335
$ws.Cells.Item(122, 11).Value = 5909.7498
$ws.Cells.Item(122, 12).Value = 16000.0005
$ws.Cells.Item(122, 13).Value = -3459.7498
$ws.Cells.Item(122, 14).Value = -20900.0005
$ws.Cells.Item(132, 8).Value = 20836364
$ws.Cells.Item(132, 9).Value = 23812202
$ws.Cells.Item(132, 11).Value = 71436606
$ws.Cells.Item(132, 13).Value = -71434076
$ws.Cells.Item(136, 8).Value = 3783.7144
$ws.Cells.Item(136, 9).Value = 1746.5
$ws.Cells.Item(136, 10).Value = 6500
$ws.Cells.Item(136, 11).Value = 5239.5
$ws.Cells.Item(136, 12).Value = 19500
$ws.Cells.Item(136, 13).Value = -2689.5
$ws.Cells.Item(136, 14).Value = -24600

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1424.84
$ws.Cells.Item(20, 9).Value = 821.4
$ws.Cells.Item(20, 10).Value = 2330
$ws.Cells.Item(20, 11).Value = 821.4
$ws.Cells.Item(20, 12).Value = 2330
$ws.Cells.Item(20, 13).Value = -574.4
$ws.Cells.Item(20, 14).Value = -2824
$ws.Cells.Item(80, 8).Value = 1169.1052
$ws.Cells.Item(80, 10).Value = 1000
$ws.Cells.Item(80, 12).Value = 1000
$ws.Cells.Item(80, 14).Value = -2996
$ws.Cells.Item(83, 8).Value = 1169.1052
$ws.Cells.Item(83, 10).Value = 1000
$ws.Cells.Item(83, 12).Value = 5000
$ws.Cells.Item(83, 14).Value = -14984

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 16131750
$ws.Cells.Item(58, 9).Value = 1534.3684
$ws.Cells.Item(58, 11).Value = 1534.3684
$ws.Cells.Item(58, 13).Value = -1331.3684
$ws.Cells.Item(122, 8).Value = 4541.222
$ws.Cells.Item(122, 9).Value = 4158
$ws.Cells.Item(122, 10).Value = 5020.25
$ws.Cells.Item(122, 11).Value = 12474
$ws.Cells.Item(122, 12).Value = 15060.75
$ws.Cells.Item(122, 13).Value = -10024
$ws.Cells.Item(122, 14).Value = -19960.75
$ws.Cells.Item(132, 8).Value = 3112.4055
$ws.Cells.Item(132, 9).Value = 1866.238
$ws.Cells.Item(132, 10).Value = 4748
$ws.Cells.Item(132, 11).Value = 5598.714
$ws.Cells.Item(132, 12).Value = 14244
$ws.Cells.Item(132, 13).Value = -3068.714
$ws.Cells.Item(132, 14).Value = -19304
$ws.Cells.Item(136, 8).Value = 16131750
$ws.Cells.Item(136, 9).Value = 1534.3684
$ws.Cells.Item(136, 11).Value = 4603.1052
$ws.Cells.Item(136, 13).Value = -2053.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 6363.3335
$ws.Cells.Item(3, 9).Value = 254
$ws.Cells.Item(3, 10).Value = 14000
$ws.Cells.Item(3, 11).Value = 762
$ws.Cells.Item(3, 12).Value = 42000
$ws.Cells.Item(3, 13).Value = -650
$ws.Cells.Item(3, 14).Value = -42224
$ws.Cells.Item(5, 8).Value = 1447.7894
$ws.Cells.Item(5, 9).Value = 514.1429000000001
$ws.Cells.Item(5, 10).Value = 4062
$ws.Cells.Item(5, 11).Value = 1542.4287
$ws.Cells.Item(5, 12).Value = 12186
$ws.Cells.Item(5, 13).Value = -1430.4287
$ws.Cells.Item(5, 14).Value = -12410
$ws.Cells.Item(86, 8).Value = 771.7143
$ws.Cells.Item(86, 9).Value = 650
$ws.Cells.Item(86, 10).Value = 820.4
$ws.Cells.Item(86, 11).Value = 1950
$ws.Cells.Item(86, 12).Value = 2461.2
$ws.Cells.Item(86, 13).Value = -764
$ws.Cells.Item(86, 14).Value = -4833.2
$ws.Cells.Item(89, 8).Value = 771.7143
$ws.Cells.Item(89, 9).Value = 650
$ws.Cells.Item(89, 10).Value = 820.4
$ws.Cells.Item(89, 11).Value = 5850
$ws.Cells.Item(89, 12).Value = 7383.599999999999
$ws.Cells.Item(89, 13).Value = 78
$ws.Cells.Item(89, 14).Value = -19239.6
$ws.Cells.Item(135, 8).Value = 1447.7894
$ws.Cells.Item(135, 9).Value = 514.1429000000001
$ws.Cells.Item(135, 10).Value = 4062
$ws.Cells.Item(135, 11).Value = 4627.2861
$ws.Cells.Item(135, 12).Value = 36558
$ws.Cells.Item(135, 13).Value = -2092.2861
$ws.Cells.Item(135, 14).Value = -41628

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3734.2222
$ws.Cells.Item(70, 9).Value = 3521.6
$ws.Cells.Item(70, 11).Value = 3521.6
$ws.Cells.Item(70, 13).Value = -3251.6
$ws.Cells.Item(73, 8).Value = 3734.2222
$ws.Cells.Item(73, 9).Value = 3521.6
$ws.Cells.Item(73, 11).Value = 3521.6
$ws.Cells.Item(73, 13).Value = -2585.6
$ws.Cells.Item(102, 8).Value = 2809.4119
$ws.Cells.Item(102, 9).Value = 2125.7144
$ws.Cells.Item(102, 10).Value = 6000
$ws.Cells.Item(102, 11).Value = 2125.7144
$ws.Cells.Item(102, 12).Value = 6000
$ws.Cells.Item(102, 13).Value = -503.7143999999998
$ws.Cells.Item(102, 14).Value = -9244
$ws.Cells.Item(132, 8).Value = 2837.074
$ws.Cells.Item(132, 9).Value = 1884.3684
$ws.Cells.Item(132, 10).Value = 5099.75
$ws.Cells.Item(132, 11).Value = 5653.1052
$ws.Cells.Item(132, 12).Value = 15299.25
$ws.Cells.Item(132, 13).Value = -3123.1052
$ws.Cells.Item(132, 14).Value = -20359.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3985
$ws.Cells.Item(40, 9).Value = 1980
$ws.Cells.Item(40, 11).Value = 1980
$ws.Cells.Item(40, 13).Value = -1844
$ws.Cells.Item(100, 8).Value = 2506.125
$ws.Cells.Item(100, 9).Value = 1111
$ws.Cells.Item(100, 10).Value = 2705.4285
$ws.Cells.Item(100, 11).Value = 1111
$ws.Cells.Item(100, 12).Value = 2705.4285
$ws.Cells.Item(100, 13).Value = -570
$ws.Cells.Item(100, 14).Value = -3787.4285
$ws.Cells.Item(132, 8).Value = 2567.2
$ws.Cells.Item(132, 9).Value = 1700.5927
$ws.Cells.Item(132, 10).Value = 3867.111
$ws.Cells.Item(132, 11).Value = 5101.7781
$ws.Cells.Item(132, 12).Value = 11601.333
$ws.Cells.Item(132, 13).Value = -2571.7781
$ws.Cells.Item(132, 14).Value = -16661.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 18009.8
$ws.Cells.Item(39, 10).Value = 18009.8
$ws.Cells.Item(39, 12).Value = 18009.8
$ws.Cells.Item(39, 14).Value = -18835.8
$ws.Cells.Item(122, 8).Value = 436717.22
$ws.Cells.Item(122, 9).Value = 501589.8
$ws.Cells.Item(122, 11).Value = 1504769.4
$ws.Cells.Item(122, 13).Value = -1502319.4
$ws.Cells.Item(132, 8).Value = 222126.28
$ws.Cells.Item(132, 9).Value = 287866.03
$ws.Cells.Item(132, 10).Value = 12954.363
$ws.Cells.Item(132, 11).Value = 863598.0900000001
$ws.Cells.Item(132, 12).Value = 38863.089
$ws.Cells.Item(132, 13).Value = -861068.0900000001
$ws.Cells.Item(132, 14).Value = -43923.089
